# Updated section 1.2 and excel file
#
# This script reproduces the edit described by the diff:
#  1. Insert a new row at row 12 for the added "2023" projection year
#     (this automatically shifts every row >= 12 down by one, and the
#     engine re-points all formula references that pointed at the old
#     rows 12-28, exactly like Excel's native "Insert Row" behaviour).
#  2. Copy the formatting (number formats / fill / styles) from the row
#     above (the previous last projected year, now row 11) down onto the
#     new row 12 so the new row keeps the same visual style as the other
#     "planned" years, then fill in the 2023 values/formulas following
#     the same pattern used by the rows above it.
#  3. Apply the other numeric/formula edits from the update:
#       - J4 (per Tier-2 totals / job slots formula constant) changed
#       - the "opportunistic" rows (now 8 and 9) get an explicit slot
#         count of 10000 instead of 0
#       - the out-year ramp (now rows 10 and 11) luminosity increments
#         drop from 80/100 to 60/60
#       - "MiniAOD replication in US" (now row 19) increases from 1.1 to 1.5
#  4. Hide the "Per Tier-2 Annual Increment" helper columns H:I.
#  5. Add the new assumption/footnote string next to "AOD Fraction on
#     disk" (now row 21) about opportunistic resources covering the
#     shortfall in 2016.
#  6. Leave the final selection on E12 (the new row), matching the
#     author's cursor position when they saved the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new row for 2023 -----------------------------------
$ws.Rows("12:12").Insert()

# --- 2. Clone formatting from the row above (2022, now row 11) --------
$ws.Range("A11:Q11").Copy()
$ws.Range("A12").PasteSpecial(-4122)   # xlPasteFormats

# --- 2b. Fill in the 2023 row, following the same formula pattern -----
$ws.Range("A12").Value = 2023
$ws.Range("B12").Value = 60
$ws.Range("C12").Formula = "=C11+B12"
$ws.Range("D12").Formula = "=D4*B12/C14"
$ws.Range("E12").Formula = "=E4*C12/C14"
$ws.Range("F12").Formula = "=(C12*C15)*C16*C17*C20*C21/C22"
$ws.Range("G12").Formula = "=(C12*C15)*C16*C17*C18*C19/C22"
$ws.Range("H12").Formula = "=J12-J11"
$ws.Range("I12").Formula = "=K12-K11"
$ws.Range("J12").Formula = "=(D12+E12)/7"
$ws.Range("K12").Formula = "=SUM(F12:G12)/7"
$ws.Range("L12").Formula = "=L11*0.9"
$ws.Range("M12").Formula = "=M11*0.9"
$ws.Range("N12").Formula = "=((H12*L12)+(I12*2*1000*M12))"
$ws.Range("O12").Formula = "=O11*0.8"
$ws.Range("P12").Formula = "=P11*0.8"
$ws.Range("Q12").Formula = "=((H12*O12)+(I12*2*1000*P12))"

# --- 3. Other value / formula edits ------------------------------------
$ws.Range("J4").Formula = "=40698/7"
$ws.Range("D8").Value = 10000
$ws.Range("D9").Value = 10000
$ws.Range("B10").Value = 60
$ws.Range("B11").Value = 60
$ws.Range("C19").Value = 1.5

# --- 4. Hide the per-Tier-2 annual increment helper columns ------------
$ws.Range("H1:I1").EntireColumn.Hidden = $true

# --- 5. New assumption footnote next to "AOD Fraction on disk" --------
$ws.Range("K21").Value = "We expect that opportunistic resources can provide some of the shortfall, especially in 2016."

# --- 6. Restore the author's final selection ---------------------------
$ws.Range("E12").Select()
